$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$used = $ws.UsedRange
$rows = $used.Rows.Count

for ($r = 1; $r -le $rows; $r++) {
    $cell = $ws.Cells.Item($r, 7)
    $v = $cell.Text
    if ($v -eq "System, dnasr281@gmail.com") {
        $cell.Value = "dnasr281@gmail.com, System"
    }
}
